$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Jugadores": add a new player row (row 3), duplicated from row 2
# (same club/scouting profile) but for a different player ("Johan Rojas2")
# and a different Transfermarket value (bug fix: wrong crest path on row 2
# is fixed for the new row, Transfermarket goes from 4 -> 30).
# ---------------------------------------------------------------------------
$wsJug = $wb.Worksheets.Item("Jugadores")

# Duplicate row 2 (values + formatting) into row 3, then patch the two
# cells that differ: A (player name) and O (Transfermarket).
$wsJug.Range("A2:BN2").Copy($wsJug.Range("A3:BN3"))
$wsJug.Range("A3").Value = "Johan Rojas2"
$wsJug.Range("O3").Value = 30

# Move the selection (cosmetic UI state captured in the saved file).
$wsJug.Range("A4").Select()

# ---------------------------------------------------------------------------
# Sheet "Entrenadores": crest-size fix. Shrink the "Fase Ofensiva" column
# (L) and wrap its header so the San Lorenzo crest / long scouting text
# shows completely.
# ---------------------------------------------------------------------------
$wsEnt = $wb.Worksheets.Item("Entrenadores")
$wsEnt.Activate()

$wsEnt.Range("L1").WrapText = $true
$wsEnt.Columns.Item(12).ColumnWidth = 8.666666666666666

# Row 15 already had a manual (custom) height because of its long wrapped
# scouting notes; narrowing column L makes that text wrap across more
# lines, so the row needs to grow to keep showing it completely.
$wsEnt.Rows.Item(15).RowHeight = 29

# Move the selection (cosmetic UI state captured in the saved file).
$wsEnt.Range("H27").Select()
